$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F (District), shifting District to G
$ws.Columns("F:F").Insert()

# New header for inserted column
$ws.Range("F2").Value = "Address"

# Populate Address values (derived from the school/address portion of column B)
$ws.Range("F3").Value = "Govt. Girls High School Humnabad"
$ws.Range("F4").Value = "G H S BelagoduSakaleshapura"
$ws.Range("F5").Value = "G H S KeragoduHolenarsipura"
$ws.Range("F6").Value = "G H S AralalusandraBidadi"
$ws.Range("F7").Value = "G H S MugalkhodMudhol"
$ws.Range("F8").Value = "G H S DhavaleshwarMudhol"
$ws.Range("F9").Value = "Govt. Adarsh Vidyalaya BemalkhedHumnabad"
$ws.Range("F10").Value = "G H S KudumalligeThirthahalli"
$ws.Range("F11").Value = "Malnad High School Hirejambur Shikaripur"
$ws.Range("F12").Value = "G H S MahalingapurMudhol"
$ws.Range("F13").Value = "G H S KakhandakiVijayapur Rural"
$ws.Range("F14").Value = "Govt. P U College(HS) BanahattiJamkhandi"
$ws.Range("F15").Value = "B Siddannaiah High School Ballupet(P)Sakaleshpur"
$ws.Range("F16").Value = "G H S GadwantiHumnabad"
$ws.Range("F17").Value = "G H S Harisandra"
$ws.Range("F18").Value = "G H S HalebelagolaChannaraya Patana"
$ws.Range("F20").Value = "G U H SchoolChannagiri"
$ws.Range("F21").Value = "G H S H Basavapura"
$ws.Range("F22").Value = "Govt. High SchoolChintakiAurad"
$ws.Range("F23").Value = "Sri Adichunchanagiri High School K KrihallyKowshika Post"
$ws.Range("F24").Value = "G U H S Extension Ramanagara Town"
$ws.Range("F25").Value = "Govt. High School KodlaSedam"
$ws.Range("F26").Value = "Shanthaveri Gopala Gowda High SchoolSydoorSagar"
$ws.Range("F27").Value = "G H S AkkimaradiMudhol"
$ws.Range("F28").Value = "G H S ThoranagatteJagalur"
$ws.Range("F29").Value = "G H S MavinakatteChannagiri"
$ws.Range("F30").Value = "Ramadurga High School"
$ws.Range("F31").Value = "Govt. Urdu High School JalwadSindgi"
$ws.Range("F32").Value = "G H S Kuntinamadu Tarikere"
$ws.Range("F33").Value = "Kanva Maharshi High School"
$ws.Range("F34").Value = "G H S RamakrishnapuraThirthahalli"
$ws.Range("F35").Value = "Sangolli Rayanna High SchoolTuppadakurahattiNavalgund"
$ws.Range("F36").Value = "S A P U C Annigeri"
$ws.Range("F38").Value = "Sanjose High School Godikoppa"
$ws.Range("F39").Value = "Jnanasindhu Rural High SchoolHallimala"
$ws.Range("F40").Value = "Govt. High School B BasapuraDavanagere North"
$ws.Range("F41").Value = "G H S KunchawaramChincholi"
$ws.Range("F42").Value = "G H S Gopanahally"
$ws.Range("F43").Value = "G H S HoranaduMoodigere"
$ws.Range("F44").Value = "Govt. High School BatgeraBasavakalyan"
$ws.Range("F45").Value = "Shree Beereshwar High School ShiradonChadachan"
$ws.Range("F46").Value = "G H S JainapurVijayapur Rural"
$ws.Range("F47").Value = "Govt. High School Hebbakodi"
$ws.Range("F48").Value = "Basaveshwar Comp. Jr. College"
$ws.Range("F49").Value = "S S H S Gulagal JambagiMudhol"
$ws.Range("F50").Value = "Govt. High School VangoorSakaleshpur"
$ws.Range("F51").Value = "N E S High School Chittapur"
$ws.Range("F52").Value = "S H High School Nainegali"
$ws.Range("F53").Value = "G H S Yarehalli"
$ws.Range("F54").Value = "Pragati H S MahalingpurMudhol"
$ws.Range("F55").Value = "G H S Sriramanagara"
$ws.Range("F56").Value = "G H S Kodiyala Karenahalli"
$ws.Range("F57").Value = "K P S ArjunagiVijayapur Rural"
$ws.Range("F58").Value = "G H S KhanadalKalaburagi South"
$ws.Range("F59").Value = "Govt. High School ShalavadiNavalgund"
$ws.Range("F60").Value = "St. Joseph High School Sakaleshapura"
$ws.Range("F61").Value = "Govt. High School (RMSA) MutturJamkhandi"
$ws.Range("F62").Value = "G U H S KerebilchiChannagiri"
